$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 425
$ws.Cells.Item(12, 10).Value = 900
$ws.Cells.Item(12, 12).Value = 900
$ws.Cells.Item(12, 14).Value = -1240
# Row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).ClearContents()
# Row 23
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()
# Row 38
$ws.Cells.Item(38, 8).Value = 1080
$ws.Cells.Item(38, 9).Value = 212.1
$ws.Cells.Item(38, 10).Value = 3249.75
$ws.Cells.Item(38, 11).Value = 636.3
$ws.Cells.Item(38, 12).Value = 9749.25
$ws.Cells.Item(38, 13).Value = -264.3
$ws.Cells.Item(38, 14).Value = -10493.25
# Row 40
$ws.Cells.Item(40, 8).Value = 4985.5884
$ws.Cells.Item(40, 10).Value = 6533.8335
$ws.Cells.Item(40, 12).Value = 6533.8335
$ws.Cells.Item(40, 14).Value = -6883.8335
# Row 55
$ws.Cells.Item(55, 8).Value = 53.285713
$ws.Cells.Item(55, 9).Value = 90
$ws.Cells.Item(55, 11).Value = 90
$ws.Cells.Item(55, 13).Value = 124
# Row 64
$ws.Cells.Item(64, 8).Value = 9218.25
$ws.Cells.Item(64, 10).Value = 9749.5
$ws.Cells.Item(64, 12).Value = 9749.5
$ws.Cells.Item(64, 14).Value = -10245.5
# Row 67
$ws.Cells.Item(67, 8).Value = 9218.25
$ws.Cells.Item(67, 10).Value = 9749.5
$ws.Cells.Item(67, 12).Value = 9749.5
$ws.Cells.Item(67, 14).Value = -11465.5
# Row 107
$ws.Cells.Item(107, 8).Value = 222.24
$ws.Cells.Item(107, 9).Value = 116.22727
$ws.Cells.Item(107, 10).Value = 999.6667
$ws.Cells.Item(107, 11).Value = 116.22727
$ws.Cells.Item(107, 12).Value = 999.6667
$ws.Cells.Item(107, 13).Value = 1803.77273
$ws.Cells.Item(107, 14).Value = -4839.6667
# Row 129
$ws.Cells.Item(129, 8).Value = 2506.111
$ws.Cells.Item(129, 9).Value = 2261.5
$ws.Cells.Item(129, 10).Value = 2995.3333
$ws.Cells.Item(129, 11).Value = 6784.5
$ws.Cells.Item(129, 12).Value = 8985.999899999999
$ws.Cells.Item(129, 13).Value = -1784.5
$ws.Cells.Item(129, 14).Value = -18985.9999
# Row 137
$ws.Cells.Item(137, 8).Value = 3238.4285
$ws.Cells.Item(137, 9).Value = 1971.8889
$ws.Cells.Item(137, 10).Value = 5518.2
$ws.Cells.Item(137, 11).Value = 5915.6667
$ws.Cells.Item(137, 12).Value = 16554.6
$ws.Cells.Item(137, 13).Value = -3365.6667
$ws.Cells.Item(137, 14).Value = -21654.6
# Row 138
$ws.Cells.Item(138, 8).Value = 2882.8462
$ws.Cells.Item(138, 10).Value = 3640
$ws.Cells.Item(138, 12).Value = 10920
$ws.Cells.Item(138, 14).Value = -21200

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 3573.75
$ws.Cells.Item(2, 9).Value = 3884.5715
$ws.Cells.Item(2, 10).Value = 1398
$ws.Cells.Item(2, 11).Value = 3884.5715
$ws.Cells.Item(2, 12).Value = 1398
$ws.Cells.Item(2, 13).Value = -3771.5715
$ws.Cells.Item(2, 14).Value = -1624
# Row 61
$ws.Cells.Item(61, 8).Value = 3255.7856
$ws.Cells.Item(61, 9).Value = 3006.037
$ws.Cells.Item(61, 11).Value = 3006.037
$ws.Cells.Item(61, 13).Value = -2794.037
# Row 102
$ws.Cells.Item(102, 8).Value = 3889.4614
$ws.Cells.Item(102, 9).Value = 2056.45
$ws.Cells.Item(102, 11).Value = 2056.45
$ws.Cells.Item(102, 13).Value = -434.4499999999998
# Row 116
$ws.Cells.Item(116, 8).Value = 3573.75
$ws.Cells.Item(116, 9).Value = 3884.5715
$ws.Cells.Item(116, 10).Value = 1398
$ws.Cells.Item(116, 11).Value = 3884.5715
$ws.Cells.Item(116, 12).Value = 1398
$ws.Cells.Item(116, 13).Value = -1590.5715
$ws.Cells.Item(116, 14).Value = -5986
# Row 136
$ws.Cells.Item(136, 8).Value = 3255.7856
$ws.Cells.Item(136, 9).Value = 3006.037
$ws.Cells.Item(136, 11).Value = 9018.110999999999
$ws.Cells.Item(136, 13).Value = -6468.110999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 3573.75
$ws.Cells.Item(3, 9).Value = 3884.5715
$ws.Cells.Item(3, 10).Value = 1398
$ws.Cells.Item(3, 11).Value = 3884.5715
$ws.Cells.Item(3, 12).Value = 1398
$ws.Cells.Item(3, 13).Value = -3770.5715
$ws.Cells.Item(3, 14).Value = -1626
# Row 94
$ws.Cells.Item(94, 8).Value = 850
$ws.Cells.Item(94, 9).Value = 850
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 850
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -399
$ws.Cells.Item(94, 14).ClearContents()
# Row 99
$ws.Cells.Item(99, 8).Value = 3952.2942
$ws.Cells.Item(99, 9).Value = 3612.6667
$ws.Cells.Item(99, 11).Value = 3612.6667
$ws.Cells.Item(99, 13).Value = -2114.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Cells.Item(107, 8).Value = 525.875
$ws.Cells.Item(107, 9).Value = 472.42856
$ws.Cells.Item(107, 10).Value = 900
$ws.Cells.Item(107, 11).Value = 472.42856
$ws.Cells.Item(107, 12).Value = 900
$ws.Cells.Item(107, 13).Value = 1447.57144
$ws.Cells.Item(107, 14).Value = -4740
# Row 119
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
# Row 132
$ws.Cells.Item(132, 8).Value = 2932.2
$ws.Cells.Item(132, 9).Value = 2427.7144
$ws.Cells.Item(132, 11).Value = 7283.1432
$ws.Cells.Item(132, 13).Value = -4753.1432
# Row 134
$ws.Cells.Item(134, 8).Value = 2927.4
$ws.Cells.Item(134, 10).Value = 5984
$ws.Cells.Item(134, 12).Value = 17952
$ws.Cells.Item(134, 14).Value = -23022
# Row 141
$ws.Cells.Item(141, 8).Value = 56649.168
$ws.Cells.Item(141, 10).Value = 56649.168
$ws.Cells.Item(141, 12).Value = 56649.168
$ws.Cells.Item(141, 14).Value = -67009.16800000001

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 1099.4286
$ws.Cells.Item(131, 10).Value = 1155.4445
$ws.Cells.Item(131, 12).Value = 3466.3335
$ws.Cells.Item(131, 14).Value = -13546.3335

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Cells.Item(97, 8).Value = 931.2857
$ws.Cells.Item(97, 9).Value = 732.0909
$ws.Cells.Item(97, 11).Value = 732.0909
$ws.Cells.Item(97, 13).Value = -236.0909
# Row 113
$ws.Cells.Item(113, 8).Value = 5585.615
$ws.Cells.Item(113, 9).Value = 1400
$ws.Cells.Item(113, 11).Value = 1400
$ws.Cells.Item(113, 13).Value = 770
# Row 132
$ws.Cells.Item(132, 8).Value = 48921.32
$ws.Cells.Item(132, 9).Value = 65502
$ws.Cells.Item(132, 10).Value = 6285.2856
$ws.Cells.Item(132, 11).Value = 196506
$ws.Cells.Item(132, 12).Value = 18855.8568
$ws.Cells.Item(132, 13).Value = -193976
$ws.Cells.Item(132, 14).Value = -23915.8568

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Cells.Item(93, 8).Value = 1331.375
$ws.Cells.Item(93, 9).Value = 965.4545000000001
$ws.Cells.Item(93, 10).Value = 1641
$ws.Cells.Item(93, 11).Value = 965.4545000000001
$ws.Cells.Item(93, 12).Value = 1641
$ws.Cells.Item(93, 13).Value = 282.5454999999999
$ws.Cells.Item(93, 14).Value = -4137
# Row 100
$ws.Cells.Item(100, 8).Value = 6749.5557
$ws.Cells.Item(100, 9).Value = 3561.5
$ws.Cells.Item(100, 10).Value = 9300
$ws.Cells.Item(100, 11).Value = 3561.5
$ws.Cells.Item(100, 12).Value = 9300
$ws.Cells.Item(100, 13).Value = -3020.5
$ws.Cells.Item(100, 14).Value = -10382
# Row 132
$ws.Cells.Item(132, 8).Value = 7428.7393
$ws.Cells.Item(132, 9).Value = 5749.375
$ws.Cells.Item(132, 11).Value = 17248.125
$ws.Cells.Item(132, 13).Value = -14718.125

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 6400.25
$ws.Cells.Item(132, 9).Value = 6400.25
$ws.Cells.Item(132, 11).Value = 19200.75
$ws.Cells.Item(132, 13).Value = -16670.75
# Row 136
$ws.Cells.Item(136, 8).Value = 8838.375
$ws.Cells.Item(136, 9).Value = 8911.25
$ws.Cells.Item(136, 11).Value = 26733.75
$ws.Cells.Item(136, 13).Value = -24183.75
